# Loan RBI, Variable Instalments
# - Inserts a new (blank) column N on the "Repayment Schedule" sheet,
#   pushing the existing "In Advance" / "Outstanding" / "Late" columns
#   one position to the right (N->O->P->Q).
# - Makes "Repayment Schedule" the active/selected sheet/tab instead of
#   "Transactions", with O6 selected on it.

$wb = $excel.ActiveWorkbook

$repayment = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (14th column).
$repayment.Columns("N").Insert()

# New column inherits the width of the column to its left (M), matching
# Excel's default behaviour when inserting a column.
$repayment.Columns("N").ColumnWidth = $repayment.Columns("M").ColumnWidth

# Make "Repayment Schedule" the active sheet/tab (previously "Transactions"
# was active), and select cell O6 on it.
$repayment.Activate() | Out-Null
$repayment.Range("O6").Select() | Out-Null
